$wb = $excel.ActiveWorkbook

# --- Sheet1 "verifyAvailablePets": browser value edge -> chrome, becomes the active/selected sheet & cell ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A2").Value = "chrome"

# --- Sheet4 "VerifySignInButtonPresence": browser value firefox -> chrome, selection moves to A11 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A2").Value = "chrome"

# --- Sheet5 "Verifylogin": remaining firefox rows -> chrome, selection moves to A6, no longer the active tab ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("A3").Value = "chrome"
$ws5.Range("A5").Value = "chrome"

# Update per-sheet selections/active cells to match the new state.
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws6 = $wb.Worksheets.Item(6)

$ws2.Range("J20").Select() | Out-Null
$ws3.Range("I25").Select() | Out-Null
$ws4.Range("A11").Select() | Out-Null
$ws5.Range("A6").Select() | Out-Null
$ws6.Range("G8").Select() | Out-Null

# Sheet1 becomes the active sheet with L22 selected (sets tabSelected on sheet1
# and clears it from the previously-active sheet5).
$ws1.Select() | Out-Null
$ws1.Range("L22").Select() | Out-Null
